$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the two text values (shared strings) that changed.
$ws.Range("A2").Value = "TettGSU"
$ws.Range("A3").Value = "TetiTLW"

# Update the selected cell to A3 (matches sheetView selection change).
$ws.Range("A3").Select()
